$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("logs")
$logs1 = $wb.Worksheets.Item("logs-1")

# ---------------------------------------------------------------------
# 1) New experiment rows (30-33) on the "logs" sheet
#    Values/formulas must be written in a specific order so that new
#    shared-string entries land at the same indices as the target file:
#      143 "Unfreeze LM head"
#      144 "lm-head"
#      145 "Unfreeze LM head, apply dropout"
#      146 "Unfreeze LM head, apply weight decay"
#      147 "0.2-0.2-0.2"
# ---------------------------------------------------------------------

# -- Row 30: Unfreeze LM head --
$logs.Range("A30").Value = "Unfreeze LM head"
$logs.Range("F30").Value = "lm-head"

# -- Row 31: Unfreeze LM head, apply dropout --
$logs.Range("A31").Value = "Unfreeze LM head, apply dropout"

# -- Row 33: Unfreeze LM head, apply weight decay (written before N31 on purpose) --
$logs.Range("A33").Value = "Unfreeze LM head, apply weight decay"

# -- Row 31/32 dropout level string --
$logs.Range("N31").Value = "0.2-0.2-0.2"

# -- Row 32 reuses the "apply dropout" note --
$logs.Range("A32").Value = "Unfreeze LM head, apply dropout"
$logs.Range("N32").Value = "0.2-0.2-0.2"

# ---- Row 30 remaining cells ----
$logs.Range("B30").Formula = '=CONCATENATE(E30,"_uf-",F30,"_ebs-",I30*K30,"_lr-",L30,"-",M30)'
$logs.Range("C30").Formula = '=CONCATENATE(E30,"_uf-",F30,"_lora-",G30,"_nepoch-",H30,"_ebs-",I30*K30,"_lr-",L30,"-",M30,"_drop-",N30,"_wd-",O30)'
$logs.Range("D30").Value = "tvtsplit"
$logs.Range("E30").Value = "bart-base"
$logs.Range("G30").Value = "none"
$logs.Range("H30").Value = 15
$logs.Range("I30").Value = 2
$logs.Range("J30").Value = 2
$logs.Range("K30").Value = 4
$logs.Range("L30").Value = 0.0004
$logs.Range("M30").Value = "step-1-0.998"
$logs.Range("N30").Value = "0.1-0.1-0.1"
$logs.Range("O30").Value = 0
$logs.Range("P30").Value = "no"

# ---- Row 31 remaining cells ----
$logs.Range("B31").Formula = '=CONCATENATE(E31,"_uf-",F31,"_ebs-",I31*K31,"_lr-",L31,"-",M31,"_drop-",N31)'
$logs.Range("C31").Formula = '=CONCATENATE(E31,"_uf-",F31,"_lora-",G31,"_nepoch-",H31,"_ebs-",I31*K31,"_lr-",L31,"-",M31,"_drop-",N31,"_wd-",O31)'
$logs.Range("D31").Value = "tvtsplit"
$logs.Range("E31").Value = "bart-base"
$logs.Range("F31").Value = "lm-head"
$logs.Range("G31").Value = "none"
$logs.Range("H31").Value = 15
$logs.Range("I31").Value = 2
$logs.Range("J31").Value = 2
$logs.Range("K31").Value = 4
$logs.Range("L31").Value = 0.0004
$logs.Range("M31").Value = "step-1-0.998"
$logs.Range("O31").Value = 0
$logs.Range("P31").Value = "no"

# ---- Row 32 remaining cells ----
$logs.Range("B32").Formula = '=CONCATENATE(E32,"_uf-",F32,"_ebs-",I32*K32,"_lr-",L32,"-",M32,"_drop-",N32,"_pt2")'
$logs.Range("C32").Formula = '=CONCATENATE(E32,"_uf-",F32,"_lora-",G32,"_nepoch-",H32,"_ebs-",I32*K32,"_lr-",L32,"-",M32,"_drop-",N32,"_wd-",O32)'
$logs.Range("D32").Value = "tvtsplit"
$logs.Range("E32").Value = "bart-base"
$logs.Range("F32").Value = "lm-head"
$logs.Range("G32").Value = "none"
$logs.Range("H32").Value = 30
$logs.Range("I32").Value = 2
$logs.Range("J32").Value = 2
$logs.Range("K32").Value = 4
$logs.Range("L32").Value = 0.0004
$logs.Range("M32").Value = "step-1-0.998"
$logs.Range("O32").Value = 0
$logs.Range("P32").Value = "no"

# ---- Row 33 remaining cells ----
$logs.Range("B33").Formula = '=CONCATENATE(E33,"_uf-",F33,"_ebs-",I33*K33,"_lr-",L33,"-",M33,"_wd-",O33)'
$logs.Range("C33").Formula = '=CONCATENATE(E33,"_uf-",F33,"_lora-",G33,"_nepoch-",H33,"_ebs-",I33*K33,"_lr-",L33,"-",M33,"_drop-",N33,"_wd-",O33)'
$logs.Range("D33").Value = "tvtsplit"
$logs.Range("E33").Value = "bart-base"
$logs.Range("F33").Value = "lm-head"
$logs.Range("G33").Value = "none"
$logs.Range("H33").Value = 15
$logs.Range("I33").Value = 2
$logs.Range("J33").Value = 2
$logs.Range("K33").Value = 4
$logs.Range("L33").Value = 0.0004
$logs.Range("M33").Value = "step-1-0.998"
$logs.Range("N33").Value = "0.1-0.1-0.1"
$logs.Range("O33").Value = 4
$logs.Range("P33").Value = "no"

# ---------------------------------------------------------------------
# 2) Re-apply the formatting that the new rows should inherit from the
#    row above them (row 29), without disturbing the values/formulas
#    that were just written.
# ---------------------------------------------------------------------
$logs.Range("A29:P29").Copy()
$logs.Range("A30:P33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$logs.Range("L30:L33").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------
# 3) Column widths on "logs" - let Excel recompute the bestFit widths
#    now that longer strings were added to columns A and B.
# ---------------------------------------------------------------------
$logs.Columns.Item(1).AutoFit()
$logs.Columns.Item(2).AutoFit()

# ---------------------------------------------------------------------
# 4) View-state changes
# ---------------------------------------------------------------------
# "logs-1" pane scroll position / selection
$logs1.Activate()
$logs1.Range("F89").Select()

# "logs" becomes the active/selected sheet (was "rouge")
$logs.Activate()
$logs.Range("B34").Select()
